# Apply updated respirometry / PvsI model-fit values to rows 2-15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1426048780487805
$ws.Range("V2").Value = 0.0002220611491829204
$ws.Range("Z2").Value = -0.1505458301548111
$ws.Range("AB2").Value = -677.947631581428
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -677.947631581428

# Row 3
$ws.Range("T3").Value = 0.1423414634146342
$ws.Range("V3").Value = 0.0002565454225970831
$ws.Range("Z3").Value = -0.1482786109104032
$ws.Range("AB3").Value = -577.9819004733592
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -577.9819004733592

# Row 4
$ws.Range("T4").Value = 0.1477268292682927
$ws.Range("V4").Value = 0.0002491214197856265
$ws.Range("Z4").Value = -0.1899555345553235
$ws.Range("AB4").Value = -762.5018142509931
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -762.5018142509931

# Row 5
$ws.Range("T5").Value = 0.1470341463414634
$ws.Range("V5").Value = 0.0002082674398172554
$ws.Range("Z5").Value = -0.1802780271119906
$ws.Range("AB5").Value = -865.6083124187626
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -865.6083124187626

# Row 6
$ws.Range("T6").Value = 0.1477560975609756
$ws.Range("V6").Value = 0.0001607142857142857
$ws.Range("Z6").Value = -0.2046258766628437
$ws.Range("AB6").Value = -1273.22767701325
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1273.22767701325

# Row 7
$ws.Range("T7").Value = 0.1455317073170732
$ws.Range("V7").Value = 0.0002247627833421192
$ws.Range("Z7").Value = -0.1161762673627402
$ws.Range("AB7").Value = -516.8839148334638
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -516.8839148334638

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("AC8").Value = "umolO2/min/m2"

# Row 9
$ws.Range("T9").Value = 0.1426048780487805
$ws.Range("V9").Value = 0.0002220611491829204
$ws.Range("Z9").Value = 0.109175566503998
$ws.Range("AB9").Value = 491.6464086838794
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 491.6464086838794

# Row 10
$ws.Range("T10").Value = 0.1423414634146342
$ws.Range("V10").Value = 0.0002565454225970831
$ws.Range("Z10").Value = 0.08224241233274605
$ws.Range("AB10").Value = 320.5764168394916
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 320.5764168394916

# Row 11
$ws.Range("T11").Value = 0.1477268292682927
$ws.Range("V11").Value = 0.0002491214197856265
$ws.Range("Z11").Value = 0.2114796645236764
$ws.Range("AB11").Value = 848.901972000876
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 848.901972000876

# Row 12
$ws.Range("T12").Value = 0.1470341463414634
$ws.Range("V12").Value = 0.0002082674398172554
$ws.Range("Z12").Value = 0.2006529113940053
$ws.Range("AB12").Value = 963.4386996357597
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 963.4386996357597

# Row 13
$ws.Range("T13").Value = 0.1477560975609756
$ws.Range("V13").Value = 0.0001607142857142857
$ws.Range("Z13").Value = 0.203544196833832
$ws.Range("AB13").Value = 1266.497224743843
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1266.497224743843

# Row 14
$ws.Range("T14").Value = 0.1455317073170732
$ws.Range("V14").Value = 0.0002247627833421192
$ws.Range("Z14").Value = 0.03620877156779433
$ws.Range("AB14").Value = 161.0977183561556
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 161.0977183561556

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.0008886321560030602
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"

# Row 8: AB8/AD8 become fully empty cells (no stored value) per source diff
$ws.Range("AB8").ClearContents()
$ws.Range("AD8").ClearContents()

Write-Output "Applied PvsI model-fit updates to rows 2-15"
